$d = $word.ActiveDocument
$q = [char]34

function Split-Runs($startPos, $parts) {
    $pos = $startPos
    $bounds = @($pos)
    foreach ($p in $parts) {
        $pos = $pos + $p.Length
        $bounds += $pos
    }
    for ($i = 0; $i -lt ($bounds.Length - 1); $i++) {
        $xr = $d.Range($bounds[$i], $bounds[$i + 1])
        $xr.Bold = 1
        $xr.Bold = 0
    }
    return $pos
}

# ---------------------------------------------------------------
# Change 1: POINT_DECL ::= "let" var "=" POINT ";"
#        -> POINT_DECL ::= "const" var "=" POINT ";"
# ---------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("POINT_DECL ::= " + $q + "let" + $q)
$start = $r.Start
$r.Text = "POINT_DECL ::= " + $q + "const" + $q

$parts1 = @(
    ("POINT_DECL ::= " + $q),
    "const",
    ($q + " var " + $q + "=" + $q + " POINT " + $q + ";" + $q)
)
Split-Runs $start $parts1 | Out-Null

# ---------------------------------------------------------------
# Change 2: BOX ::= "box (" POINT ", " POINT ");"
#        -> BOX ::= "box" "(" POINT ", " POINT ");"
# ---------------------------------------------------------------
$paras = $d.Paragraphs
$n = $paras.Count
$boxPara = $null
for ($i = 1; $i -le $n; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "BOX ::=*") { $boxPara = $p; break }
}
$boxStart = $boxPara.Range.Start

$fr = $boxPara.Range
$fr.Find.Execute(" (")
$fr.Text = $q + " " + $q + "("

$parts2 = @(
    "BOX ::= ",
    ($q + "box"),
    $q,
    " ",
    $q,
    "(",
    $q,
    (" POINT " + $q + ", " + $q + " POINT " + $q + ")"),
    ";",
    $q,
    " "
)
Split-Runs $boxStart $parts2 | Out-Null

# ---------------------------------------------------------------
# Change 3: ROUND ::= "roundAbout" string ...
#        -> ROUND ::= "roundabout" string ...
# ---------------------------------------------------------------
$eps = [char]0x03B5
$r3 = $d.Content
$r3.Find.Execute("ROUND ::= " + $q + "roundAbout" + $q)
$start3 = $r3.Start
$r3.Text = "ROUND ::= " + $q + "roundabout" + $q

$parts3 = @(
    ("ROUND ::= " + $q + "round"),
    "a",
    ("bout" + $q + " string"),
    (" " + $q + "{" + $q + " CIRC " + $q + "};" + $q + " "),
    "ROUND ",
    ("| " + $eps)
)
Split-Runs $start3 $parts3 | Out-Null

# ---------------------------------------------------------------
# Change 4: CIRC ::= "circ (" POINT "," EXPR ");"
#        -> CIRC ::= "circ" "(" POINT "," EXPR ");"
# ---------------------------------------------------------------
$circPara = $null
$paras = $d.Paragraphs
$n = $paras.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "CIRC ::=*") { $circPara = $p; break }
}
$circStart = $circPara.Range.Start

$fr4 = $circPara.Range
$fr4.Find.Execute("circ (")
$fr4.Text = "circ" + $q + " " + $q + "("

$parts4 = @(
    "CIRC ::= ",
    ($q + "circ"),
    $q,
    " ",
    $q,
    ("(" + $q + " POINT "),
    ($q + "," + $q + " "),
    "EXPR",
    " ",
    ($q + ");" + $q)
)
Split-Runs $circStart $parts4 | Out-Null

# ---------------------------------------------------------------
# Change 5: var = {A,...,Z,a,...,z}+{0,...,9}*
#        -> var = _{[A-Za-z]+[0-9]*}_
# ---------------------------------------------------------------
$varPara = $null
$paras = $d.Paragraphs
$n = $paras.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "var = {A*") { $varPara = $p; break }
}
$varStart = $varPara.Range.Start
$varLen = $varPara.Range.Text.Length
$r5 = $d.Range($varStart, $varStart + $varLen)
$r5.Text = "var = _{[A-Za-z]+[0-9]*}_"

$parts5 = @("var = ", "_{[A-Za-z]+[0-9]*}_")
Split-Runs $varStart $parts5 | Out-Null

Write-Output "done"
